# precreate masks and train with 224x224
# Updates the Scores sheet:
#  - A2: append "(288x288)" on a new line (datasubset run)
#  - B2: trailing newline added to the train-log text
#  - A3: append "(288x288)" as a separately-formatted run (rich text)
#  - B3: trailing newline added to the train-log text
#  - New row 4: full-train-data run at 256x256 (A4 rich text, B4 log)
#  - New row 5: full-train-data run at 224x224 (A5 bold, B5 log)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: datasubset model, now trained at 288x288 ----
$ws.Range("A2").Value = "Unet Efficientnet-b1 `nTrained on datasubset`n(288x288)"

$ws.Range("B2").Value = "Epoch 1 | Train Loss: 0.452 | Valid Loss: 0.271 | Combined metric: 0.574 | Dice: 0.179 (LB 0.505, SB 0.236, S 0.236) | Hausdorff: 0.162 (LB 0.153, SB 0.095, S 0.239)`nEpoch 2 | Train Loss: 0.167 | Valid Loss: 0.280 | Combined metric: 0.540 | Dice: 0.205 (LB 0.499, SB 0.285, S 0.147) | Hausdorff: 0.237 (LB 0.194, SB 0.160, S 0.357)`nEpoch 3 | Train Loss: 0.114 | Valid Loss: 0.181 | Combined metric: 0.616 | Dice: 0.241 (LB 0.622, SB 0.207, S 0.561) | Hausdorff: 0.135 (LB 0.174, SB 0.099, S 0.131)`nEpoch 4 | Train Loss: 0.091 | Valid Loss: 0.167 | Combined metric: 0.636 | Dice: 0.255 (LB 0.648, SB 0.265, S 0.578) | Hausdorff: 0.110 (LB 0.126, SB 0.132, S 0.072)`nEpoch 5 | Train Loss: 0.079 | Valid Loss: 0.170 | Combined metric: 0.646 | Dice: 0.253 (LB 0.649, SB 0.237, S 0.605) | Hausdorff: 0.091 (LB 0.083, SB 0.089, S 0.101)`n"

# ---- Row 3: full-train-data model, now labelled with its 288x288 resolution ----
# The resolution suffix is a separate run so it can carry its own (un-bolded) font.
$a3run1 = "Unet Efficientnet-b1 `nTrained on full train data`n"
$a3run2 = "(288x288)"
$ws.Range("A3").Value = $a3run1 + $a3run2
$a3c2 = $ws.Range("A3").Characters(($a3run1.Length + 1), $a3run2.Length)
$a3c2.Font.Name = "Arial"
$a3c2.Font.Size = 10

$ws.Range("B3").Value = "Epoch 1 | Train Loss: 0.203 | Valid Loss: 0.123 | Combined metric: 0.654 | Dice: 0.292 (LB 0.675, SB 0.506, S 0.660) | Hausdorff: 0.104 (LB 0.093, SB 0.164, S 0.054)`nEpoch 2 | Train Loss: 0.108 | Valid Loss: 0.114 | Combined metric: 0.671 | Dice: 0.304 (LB 0.688, SB 0.566, S 0.693) | Hausdorff: 0.084 (LB 0.108, SB 0.073, S 0.071)`nEpoch 3 | Train Loss: 0.096 | Valid Loss: 0.110 | Combined metric: 0.653 | Dice: 0.304 (LB 0.612, SB 0.565, S 0.694) | Hausdorff: 0.115 (LB 0.214, SB 0.076, S 0.054)`nEpoch 4 | Train Loss: 0.087 | Valid Loss: 0.115 | Combined metric: 0.668 | Dice: 0.300 (LB 0.687, SB 0.524, S 0.670) | Hausdorff: 0.087 (LB 0.109, SB 0.096, S 0.057)`nEpoch 5 | Train Loss: 0.081 | Valid Loss: 0.107 | Combined metric: 0.689 | Dice: 0.311 (LB 0.703, SB 0.581, S 0.728) | Hausdorff: 0.060 (LB 0.068, SB 0.065, S 0.047)`n"

# ---- Row 4 (new): full-train-data model at 256x256 ----
$a4run1 = "Unet Efficientnet-b1 `nTrained on full train data`n"
$a4run2 = "(256x256)"
$ws.Range("A4").Value = $a4run1 + $a4run2
$a4c2 = $ws.Range("A4").Characters(($a4run1.Length + 1), $a4run2.Length)
$a4c2.Font.Name = "Arial"
$a4c2.Font.Size = 10

$ws.Range("B4").Value = "Epoch 1 | Train Loss: 0.211 | Valid Loss: 0.139 | Combined metric: 0.658 | Dice: 0.291 (LB 0.588, SB 0.433, S 0.600) | Hausdorff: 0.097 (LB 0.105, SB 0.111, S 0.075)`nEpoch 2 | Train Loss: 0.110 | Valid Loss: 0.129 | Combined metric: 0.657 | Dice: 0.285 (LB 0.650, SB 0.486, S 0.654) | Hausdorff: 0.095 (LB 0.130, SB 0.097, S 0.056)`nEpoch 3 | Train Loss: 0.098 | Valid Loss: 0.124 | Combined metric: 0.609 | Dice: 0.297 (LB 0.583, SB 0.393, S 0.611) | Hausdorff: 0.183 (LB 0.217, SB 0.231, S 0.100)`nEpoch 4 | Train Loss: 0.089 | Valid Loss: 0.122 | Combined metric: 0.677 | Dice: 0.297 (LB 0.687, SB 0.558, S 0.651) | Hausdorff: 0.070 (LB 0.068, SB 0.069, S 0.071)`nEpoch 5 | Train Loss: 0.083 | Valid Loss: 0.109 | Combined metric: 0.668 | Dice: 0.306 (LB 0.686, SB 0.569, S 0.685) | Hausdorff: 0.090 (LB 0.097, SB 0.119, S 0.055)`n"

$ws.Range("C4").Value = 0.668
$ws.Range("D4").Value = 0.8153
$ws.Range("E4").Value = 0.80775

# ---- Row 5 (new): full-train-data model at 224x224 (this commit's change) ----
$ws.Range("A5").Value = "Unet Efficientnet-b1 `nTrained on full train data`n(224x224)"
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10

$ws.Range("B5").Value = "Epoch 1 | Train Loss: 0.210 | Valid Loss: 0.125 | Combined metric: 0.662 | Dice: 0.295 (LB 0.659, SB 0.527, S 0.623) | Hausdorff: 0.093 (LB 0.140, SB 0.072, S 0.068)`nEpoch 2 | Train Loss: 0.111 | Valid Loss: 0.126 | Combined metric: 0.666 | Dice: 0.294 (LB 0.663, SB 0.545, S 0.615) | Hausdorff: 0.085 (LB 0.076, SB 0.115, S 0.065)`nEpoch 3 | Train Loss: 0.098 | Valid Loss: 0.115 | Combined metric: 0.682 | Dice: 0.303 (LB 0.691, SB 0.565, S 0.690) | Hausdorff: 0.066 (LB 0.078, SB 0.070, S 0.050)`nEpoch 4 | Train Loss: 0.088 | Valid Loss: 0.127 | Combined metric: 0.660 | Dice: 0.282 (LB 0.660, SB 0.520, S 0.659) | Hausdorff: 0.088 (LB 0.060, SB 0.136, S 0.068)`nEpoch 5 | Train Loss: 0.085 | Valid Loss: 0.109 | Combined metric: 0.682 | Dice: 0.315 (LB 0.687, SB 0.564, S 0.699) | Hausdorff: 0.073 (LB 0.075, SB 0.102, S 0.041)`n"

$ws.Range("C5").Value = 0.682
$ws.Range("D5").Value = 0.82298
$ws.Range("E5").Value = 0.81752

# ---- Row heights / wrap for the (now taller, 2-line suffix) description column ----
$ws.Range("C1").WrapText = $true
$ws.Range("A2:B5").WrapText = $true
$ws.Rows.Item(2).RowHeight = 68.65
$ws.Rows.Item(3).RowHeight = 68.65
$ws.Rows.Item(4).RowHeight = 68.65
$ws.Rows.Item(5).RowHeight = 68.65

# ---- Selection moves to E11 in the saved file ----
$ws.Range("E11").Select()

Write-Output "edit applied"
